$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows at position 4 (for GoogleCloud, Octoparse, Appsheet entries)
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# Insert 1 new row at position 11 (for the newly added "限定公開" entry)
$ws.Rows.Item(11).Insert()

# Remove existing hyperlinks on the sheet (deleting via a range clears all; we re-add them all below)
$ws.Range("F2").Hyperlinks.Delete()

# Row 2
$ws.Cells.Item(2,1).Value = "2026-01-05 12:41:53"
$ws.Cells.Item(2,2).Value = "AIを活用した社内備品管理アプリ開発の相談"
$ws.Cells.Item(2,3).Value = "システム開発"
$ws.Cells.Item(2,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(2,5).Value = "期限情報なし"
$ws.Cells.Item(2,6).Value = "https://www.lancers.jp/work/detail/5465005"
$ws.Cells.Item(2,7).Value = 388
$ws.Cells.Item(2,8).Value = "🔥AI,Ai ◆開発 ◇アプリ"

# Row 3
$ws.Cells.Item(3,1).Value = "2026-01-05 12:41:53"
$ws.Cells.Item(3,2).Value = "【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています"
$ws.Cells.Item(3,3).Value = "システム開発"
$ws.Cells.Item(3,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(3,5).Value = "期限情報なし"
$ws.Cells.Item(3,6).Value = "https://www.lancers.jp/work/detail/5405023"
$ws.Cells.Item(3,7).Value = 178
$ws.Cells.Item(3,8).Value = "★bot ◆ツール"

# Row 4
$ws.Cells.Item(4,1).Value = "2026-01-05 12:41:53"
$ws.Cells.Item(4,2).Value = "GoogleCloudを利用したアジャイル開発共通基盤のSREエンジニアの募集"
$ws.Cells.Item(4,3).Value = "システム開発"
$ws.Cells.Item(4,4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(4,5).Value = "期限情報なし"
$ws.Cells.Item(4,6).Value = "https://www.lancers.jp/work/detail/5457458"
$ws.Cells.Item(4,7).Value = 75
$ws.Cells.Item(4,8).Value = "◆開発"

# Row 5
$ws.Cells.Item(5,1).Value = "2026-01-05 12:41:53"
$ws.Cells.Item(5,2).Value = "Octoparseを使ったスクレイピングシステムの構築"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,6).Value = "https://www.lancers.jp/work/detail/5465301"
$ws.Cells.Item(5,7).Value = 58
$ws.Cells.Item(5,8).Value = "◆スクレイピング"

# Row 6
$ws.Cells.Item(6,1).Value = "2026-01-05 12:41:53"
$ws.Cells.Item(6,2).Value = "【急募】Appsheetで見積もりアプリを作成してくれる方"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,6).Value = "https://www.lancers.jp/work/detail/5465442"
$ws.Cells.Item(6,7).Value = 30
$ws.Cells.Item(6,8).Value = "◇アプリ"

# Row 7
$ws.Cells.Item(7,1).Value = "2026-01-05 12:41:53"
$ws.Cells.Item(7,2).Value = "初回 bubbleで構築したサイトの修正対応"
$ws.Cells.Item(7,3).Value = "システム開発"
$ws.Cells.Item(7,4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(7,5).Value = "期限情報なし"
$ws.Cells.Item(7,6).Value = "https://www.lancers.jp/work/detail/5465187"
$ws.Cells.Item(7,7).Value = 30
$ws.Cells.Item(7,8).Value = "◇サイト"

# Row 8
$ws.Cells.Item(8,1).Value = "2026-01-05 12:41:53"
$ws.Cells.Item(8,2).Value = "【急募】メール問い合わせ時の自動SMS送信システム構築"
$ws.Cells.Item(8,3).Value = "システム開発"
$ws.Cells.Item(8,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(8,5).Value = "期限情報なし"
$ws.Cells.Item(8,6).Value = "https://www.lancers.jp/work/detail/5464796"
$ws.Cells.Item(8,7).Value = 33
$ws.Cells.Item(8,8).ClearContents()

# Row 9
$ws.Cells.Item(9,1).Value = "2026-01-05 12:41:53"
$ws.Cells.Item(9,2).Value = "金融機関の入出金伝票印刷システム構築依頼"
$ws.Cells.Item(9,3).Value = "システム開発"
$ws.Cells.Item(9,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(9,5).Value = "期限情報なし"
$ws.Cells.Item(9,6).Value = "https://www.lancers.jp/work/detail/5464833"
$ws.Cells.Item(9,7).Value = 28
$ws.Cells.Item(9,8).ClearContents()

# Row 10
$ws.Cells.Item(10,1).Value = "2026-01-05 12:41:53"
$ws.Cells.Item(10,2).Value = "【急募】クラウドウェア内製化推進のための技術サポート依頼"
$ws.Cells.Item(10,3).Value = "システム開発"
$ws.Cells.Item(10,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(10,5).Value = "期限情報なし"
$ws.Cells.Item(10,6).Value = "https://www.lancers.jp/work/detail/5465210"
$ws.Cells.Item(10,7).Value = 25
$ws.Cells.Item(10,8).ClearContents()

# Row 11
$ws.Cells.Item(11,1).Value = "2026-01-05 12:41:53"
$ws.Cells.Item(11,2).Value = "限定公開 限定公開の仕事"
$ws.Cells.Item(11,3).Value = "システム開発"
$ws.Cells.Item(11,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(11,5).Value = "期限情報なし"
$ws.Cells.Item(11,6).Value = "https://www.lancers.jp/work/detail/5465372"
$ws.Cells.Item(11,7).Value = 18
$ws.Cells.Item(11,8).ClearContents()

# Row 12
$ws.Cells.Item(12,1).Value = "2026-01-05 12:41:53"
$ws.Cells.Item(12,2).Value = "RobloxアクションRPG制作(MVP/完成版前提プロジェクト)"
$ws.Cells.Item(12,3).Value = "システム開発"
$ws.Cells.Item(12,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(12,5).Value = "期限情報なし"
$ws.Cells.Item(12,6).Value = "https://www.lancers.jp/work/detail/5465063"
$ws.Cells.Item(12,7).Value = 18
$ws.Cells.Item(12,8).ClearContents()

# Row 13
$ws.Cells.Item(13,1).Value = "2026-01-05 12:41:53"
$ws.Cells.Item(13,2).Value = "【準委任】音声データ収集プロジェクトのPM・ディレクター募集"
$ws.Cells.Item(13,3).Value = "システム開発"
$ws.Cells.Item(13,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(13,5).Value = "期限情報なし"
$ws.Cells.Item(13,6).Value = "https://www.lancers.jp/work/detail/5465028"
$ws.Cells.Item(13,7).Value = 18
$ws.Cells.Item(13,8).ClearContents()

# Re-add hyperlinks for F2:F13 in order, then reapply the Hyperlink style so the
# resulting cell style index matches the workbook's existing "Hyperlink" cell style.
$ws.Hyperlinks.Add($ws.Cells.Item(2,6), "https://www.lancers.jp/work/detail/5465005")
$ws.Hyperlinks.Add($ws.Cells.Item(3,6), "https://www.lancers.jp/work/detail/5405023")
$ws.Hyperlinks.Add($ws.Cells.Item(4,6), "https://www.lancers.jp/work/detail/5457458")
$ws.Hyperlinks.Add($ws.Cells.Item(5,6), "https://www.lancers.jp/work/detail/5465301")
$ws.Hyperlinks.Add($ws.Cells.Item(6,6), "https://www.lancers.jp/work/detail/5465442")
$ws.Hyperlinks.Add($ws.Cells.Item(7,6), "https://www.lancers.jp/work/detail/5465187")
$ws.Hyperlinks.Add($ws.Cells.Item(8,6), "https://www.lancers.jp/work/detail/5464796")
$ws.Hyperlinks.Add($ws.Cells.Item(9,6), "https://www.lancers.jp/work/detail/5464833")
$ws.Hyperlinks.Add($ws.Cells.Item(10,6), "https://www.lancers.jp/work/detail/5465210")
$ws.Hyperlinks.Add($ws.Cells.Item(11,6), "https://www.lancers.jp/work/detail/5465372")
$ws.Hyperlinks.Add($ws.Cells.Item(12,6), "https://www.lancers.jp/work/detail/5465063")
$ws.Hyperlinks.Add($ws.Cells.Item(13,6), "https://www.lancers.jp/work/detail/5465028")

$ws.Cells.Item(2,6).Style = "Hyperlink"
$ws.Cells.Item(3,6).Style = "Hyperlink"
$ws.Cells.Item(4,6).Style = "Hyperlink"
$ws.Cells.Item(5,6).Style = "Hyperlink"
$ws.Cells.Item(6,6).Style = "Hyperlink"
$ws.Cells.Item(7,6).Style = "Hyperlink"
$ws.Cells.Item(8,6).Style = "Hyperlink"
$ws.Cells.Item(9,6).Style = "Hyperlink"
$ws.Cells.Item(10,6).Style = "Hyperlink"
$ws.Cells.Item(11,6).Style = "Hyperlink"
$ws.Cells.Item(12,6).Style = "Hyperlink"
$ws.Cells.Item(13,6).Style = "Hyperlink"

# Column width adjustments (stored width = ColumnWidth + 0.8333333333333333)
$ws.Columns.Item(2).ColumnWidth = 40.166666666666664
$ws.Columns.Item(4).ColumnWidth = 29.166666666666668

Write-Output "done"